$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0334433051472082
$ws.Range("C2").Value = 0.0334433051472082
$ws.Range("D2").Value = 2.0903769276243
$ws.Range("F2").Value = 0.1472

$ws.Range("B3").Value = 2.33580962694684
$ws.Range("C3").Value = 0.0159986960749783
